# Applies the StructureDefinition-average-wholesale-price.xlsx update:
#   - URL changes from ibm.com to linuxforhealth.org
#   - Version bumps from 7.0.0 to 8.0.0
#   - Date changes to 2022-11-10T16:00:46+00:00
#   - Publisher changes from "Alvearie Team" to "LinuxForHealth Team"
#   - The Elements sheet's "Extension" row Constraint(s) cell is cleared

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/average-wholesale-price"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

$elements = $wb.Worksheets.Item("Elements")
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/average-wholesale-price"
$elements.Range("AI2").Value = ""
